$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("D1").Value = "XML File Fields Transformed"

# --- Row 2: OrderID (Match) ---
$ws.Range("C2").Value = 2002
$ws.Range("D2").Value = 2002
$ws.Range("E2").Value = "Match"

# --- Row 3: Status (Match) ---
$ws.Range("C3").Value = "PROCESSING"
$ws.Range("D3").Value = "PROCESSING"
$ws.Range("E3").Value = "Match"

# --- Row 4: CustomerID (Mismatch) ---
$ws.Range("D4").Value = 98765
$ws.Range("E4").Value = "Mismatch"

# --- Row 5: CustomerName (Match) ---
$ws.Range("C5").Value = "Alice Smith"
$ws.Range("D5").Value = "Alice Smith"
$ws.Range("E5").Value = "Match"

# --- Row 6: Email (Match) ---
$ws.Range("C6").Value = "alice.smith@example.com"
$ws.Range("D6").Value = "alice.smith@example.com"
$ws.Range("E6").Value = "Match"

# --- Row 7: Phone (Match) ---
$ws.Range("C7").Value = 15556789
$ws.Range("D7").Value = 15556789
$ws.Range("E7").Value = "Match"

# --- Row 8: OrderDate (Mismatch) ---
$ws.Range("D8").Value = "15-03-2025"
$ws.Range("E8").Value = "Mismatch"

# --- Row 9: TotalAmount (Match) ---
$ws.Range("C9").Value = 500
$ws.Range("D9").Value = 500
$ws.Range("E9").Value = "Match"

# --- Row 10: Currency (Match) ---
$ws.Range("D10").Value = "USD"
$ws.Range("E10").Value = "Match"

# --- Row 11: ShippingMethod (Mismatch) ---
$ws.Range("D11").Value = "Standard"
$ws.Range("E11").Value = "Mismatch"

# --- Row 12: TrackingNumber (Mismatch) ---
$ws.Range("D12").Value = "STD987654321"
$ws.Range("E12").Value = "Mismatch"

# --- Row 13: Street (Mismatch) ---
$ws.Range("D13").Value = "456 Oak Street"
$ws.Range("E13").Value = "Mismatch"

# --- Row 14: City (Match) ---
$ws.Range("C14").Value = "Los Angeles"
$ws.Range("D14").Value = "Los Angeles"
$ws.Range("E14").Value = "Match"

# --- Row 15: State (Match) ---
$ws.Range("C15").Value = "CA"
$ws.Range("D15").Value = "CA"
$ws.Range("E15").Value = "Match"

# --- Row 16: Zip (Match) ---
$ws.Range("C16").Value = 90001
$ws.Range("D16").Value = 90001
$ws.Range("E16").Value = "Match"

# --- Row 17: Country (Match) ---
$ws.Range("D17").Value = "USA"
$ws.Range("E17").Value = "Match"

# --- Row 18: LineNumber (Match) ---
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = "Match"

# --- Row 19: ProductID (Match) ---
$ws.Range("C19").Value = 55555
$ws.Range("D19").Value = 55555
$ws.Range("E19").Value = "Match"

# --- Row 20: ProductName (Mismatch) ---
$ws.Range("C20").Value = "BlutoothSpeaker"
$ws.Range("D20").Value = "Bluetooth Speaker"
$ws.Range("E20").Value = "Mismatch"

# --- Row 21: Category (Match) ---
$ws.Range("C21").Value = "Audio"
$ws.Range("D21").Value = "Audio"
$ws.Range("E21").Value = "Match"

# --- Row 22: Quantity (Mismatch) - D22 stored as text "1" ---
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1"
$ws.Range("E22").Value = "Mismatch"

# --- Row 23: UnitPrice (Match) ---
$ws.Range("C23").Value = 100
$ws.Range("D23").Value = 100
$ws.Range("E23").Value = "Match"

# --- Row 24: Currency (Match) ---
$ws.Range("D24").Value = "USD"
$ws.Range("E24").Value = "Match"

# --- Row 25: DiscountAmount (Mismatch) - D25 stored as text "0.00" ---
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.00"
$ws.Range("E25").Value = "Mismatch"

# --- Row 26: DiscountPercentage (Mismatch) - D26 stored as text "0" ---
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("E26").Value = "Mismatch"
